$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 265 (pushes existing rows 265:326 down to 266:327),
# adding this week's new price record for Betarraga at Macroferia Regional de Talca.
$ws.Rows.Item(265).Insert()

# Populate the newly inserted row with the new weekly record
$ws.Range("A265").Value = 5
$ws.Range("B265").Value = 'Macroferia Regional de Talca'
$ws.Range("C265").Value = 'Maule'
$ws.Range("D265").Value = 44711
$ws.Range("E265").Value = 7
$ws.Range("F265").Value = 100114014
$ws.Range("G265").Value = 'Betarraga'
$ws.Range("H265").Value = 'Sin especificar'
$ws.Range("I265").Value = 'Primera'
$ws.Range("J265").Value = 4000
$ws.Range("K265").Value = 650
$ws.Range("L265").Value = 650
$ws.Range("M265").Value = 650
$ws.Range("N265").Value = '$/paquete 5 unidades'
$ws.Range("O265").Value = 'Región del Maule'
$ws.Range("P265").Value = 130
$ws.Range("Q265").Value = 5
$ws.Range("R265").Value = 'Hortaliza'
